$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D28").Value = "[논문 리뷰] Distributed neighbor selection in multi agent network"
$ws.Range("E28").Value = "https://ropiens.tistory.com/267"

$ws.Range("D37").Value = "[Paper Review] HDMixer: Hierarchical Dependency with Extendable Patch for Multivariate Time Series Forecasting"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?mod=document&uid=3157"
